$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.3846535682678223
$ws.Range("E2").Value = 521.1688494529626
$ws.Range("F2").Value = 0.01726134875182288
$ws.Range("G2").Value = 0.01467704100518867
$ws.Range("H2").Value = 0.01384062031858268
$ws.Range("I2").Value = 0.01280939414495729
$ws.Range("J2").Value = 0.01228255202089138
$ws.Range("K2").Value = 0.01183247548436967
$ws.Range("L2").Value = 0.01123030768513341
$ws.Range("M2").Value = 0.01123030768513341
$ws.Range("N2").Value = 0.01092597154557071
$ws.Range("O2").Value = 0.01091519473500512
$ws.Range("P2").Value = 0.01082650484581931
$ws.Range("Q2").Value = 0.01046071477024502
$ws.Range("R2").Value = 0.01046071477024502
$ws.Range("S2").Value = 0.01037973351206873
$ws.Range("T2").Value = 0.01037196913675044
$ws.Range("U2").Value = 0.01019320927692346
$ws.Range("V2").Value = 0.01019320927692346
$ws.Range("W2").Value = 0.01019320927692346
$ws.Range("X2").Value = 0.01015923683144176
$ws.Range("Y2").Value = 0.01015923683144176

# Row 3
$ws.Range("C3").Value = 0.3593742847442627
$ws.Range("E3").Value = 524.5924523845588
$ws.Range("F3").Value = 0.016585056546074
$ws.Range("G3").Value = 0.01478434532852796
$ws.Range("H3").Value = 0.01314350550964602
$ws.Range("I3").Value = 0.01268707457961382
$ws.Range("J3").Value = 0.01177644641711992
$ws.Range("K3").Value = 0.01172341303856955
$ws.Range("L3").Value = 0.01141759694140061
$ws.Range("M3").Value = 0.01131614967141499
$ws.Range("N3").Value = 0.01114894184153122
$ws.Range("O3").Value = 0.01104190777532934
$ws.Range("P3").Value = 0.01082862586967238
$ws.Range("Q3").Value = 0.01074449385412301
$ws.Range("R3").Value = 0.01062864311849462
$ws.Range("S3").Value = 0.01044674987314853
$ws.Range("T3").Value = 0.01038995545239249
$ws.Range("U3").Value = 0.01038995545239249
$ws.Range("V3").Value = 0.01035682077460945
$ws.Range("W3").Value = 0.01028991148812799
$ws.Range("X3").Value = 0.01026849241849766
$ws.Range("Y3").Value = 0.01022597373069315

# Row 4
$ws.Range("C4").Value = 0.4062278270721436
$ws.Range("E4").Value = 527.7575298668398
$ws.Range("F4").Value = 0.01661958254890368
$ws.Range("G4").Value = 0.0139121335172731
$ws.Range("H4").Value = 0.01244567316870655
$ws.Range("I4").Value = 0.01244567316870655
$ws.Range("J4").Value = 0.01186942887502129
$ws.Range("K4").Value = 0.01127658878860997
$ws.Range("L4").Value = 0.01106090081945583
$ws.Range("M4").Value = 0.01103715200955758
$ws.Range("N4").Value = 0.01091699760276501
$ws.Range("O4").Value = 0.01091699760276501
$ws.Range("P4").Value = 0.01069368041022797
$ws.Range("Q4").Value = 0.01056073968213496
$ws.Range("R4").Value = 0.01046407028316742
$ws.Range("S4").Value = 0.01043113264159436
$ws.Range("T4").Value = 0.01042692078635604
$ws.Range("U4").Value = 0.01037934891977649
$ws.Range("V4").Value = 0.01031614302478213
$ws.Range("W4").Value = 0.01031614302478213
$ws.Range("X4").Value = 0.01031614302478213
$ws.Range("Y4").Value = 0.01028767114750175

# Row 5
$ws.Range("C5").Value = 0.4626882076263428
$ws.Range("E5").Value = 504.7593269687386
$ws.Range("F5").Value = 0.01707601669249577
$ws.Range("G5").Value = 0.01328809265566422
$ws.Range("H5").Value = 0.01253726973363074
$ws.Range("I5").Value = 0.01225157645516308
$ws.Range("J5").Value = 0.01206533397682729
$ws.Range("K5").Value = 0.01141330369567833
$ws.Range("L5").Value = 0.01079701333857103
$ws.Range("M5").Value = 0.01079701333857103
$ws.Range("N5").Value = 0.01074181164407981
$ws.Range("O5").Value = 0.01058120757559798
$ws.Range("P5").Value = 0.01049853691585866
$ws.Range("Q5").Value = 0.01008521013865547
$ws.Range("R5").Value = 0.01008521013865547
$ws.Range("S5").Value = 0.009992220940207576
$ws.Range("T5").Value = 0.009992220940207576
$ws.Range("U5").Value = 0.009924650545742326
$ws.Range("V5").Value = 0.00989873450909152
$ws.Range("W5").Value = 0.009879645654146001
$ws.Range("X5").Value = 0.009849121767458328
$ws.Range("Y5").Value = 0.009839363098805818

# Row 6
$ws.Range("C6").Value = 0.3774311542510986
$ws.Range("E6").Value = 519.4926446873378
$ws.Range("F6").Value = 0.01678409413846915
$ws.Range("G6").Value = 0.01489863120862301
$ws.Range("H6").Value = 0.01392785216993685
$ws.Range("I6").Value = 0.01311325255557173
$ws.Range("J6").Value = 0.01204080064002014
$ws.Range("K6").Value = 0.01170064599677505
$ws.Range("L6").Value = 0.01153513961602028
$ws.Range("M6").Value = 0.01113532842885985
$ws.Range("N6").Value = 0.01084072855808763
$ws.Range("O6").Value = 0.01077604366151271
$ws.Range("P6").Value = 0.01058345459692351
$ws.Range("Q6").Value = 0.01058345459692351
$ws.Range("R6").Value = 0.01041273485410035
$ws.Range("S6").Value = 0.01041273485410035
$ws.Range("T6").Value = 0.01030273163960288
$ws.Range("U6").Value = 0.01023549666323172
$ws.Range("V6").Value = 0.01023549666323172
$ws.Range("W6").Value = 0.01019346426386551
$ws.Range("X6").Value = 0.01015082083695588
$ws.Range("Y6").Value = 0.01012656227460697

# Row 7
$ws.Range("C7").Value = 0.3906140327453613
$ws.Range("E7").Value = 538.6989914838796
$ws.Range("F7").Value = 0.01721052787329705
$ws.Range("G7").Value = 0.01457097567962012
$ws.Range("H7").Value = 0.01328348601365651
$ws.Range("I7").Value = 0.01266941080418793
$ws.Range("J7").Value = 0.01228864105391936
$ws.Range("K7").Value = 0.0119086051992593
$ws.Range("L7").Value = 0.01127716944474436
$ws.Range("M7").Value = 0.01127220351041182
$ws.Range("N7").Value = 0.01127220351041182
$ws.Range("O7").Value = 0.01115788806091579
$ws.Range("P7").Value = 0.01115212997475195
$ws.Range("Q7").Value = 0.01109014657075157
$ws.Range("R7").Value = 0.01083349829690484
$ws.Range("S7").Value = 0.01083349829690484
$ws.Range("T7").Value = 0.01063707436966335
$ws.Range("U7").Value = 0.01063707436966335
$ws.Range("V7").Value = 0.01060425524046618
$ws.Range("W7").Value = 0.01055295427649654
$ws.Range("X7").Value = 0.01051305923166879
$ws.Range("Y7").Value = 0.01050095499968576

# Row 8
$ws.Range("C8").Value = 0.3749752044677734
$ws.Range("E8").Value = 516.766452802427
$ws.Range("F8").Value = 0.01705956309028638
$ws.Range("G8").Value = 0.01438052092420931
$ws.Range("H8").Value = 0.01326188934528946
$ws.Range("I8").Value = 0.01191432837915967
$ws.Range("J8").Value = 0.01191432837915967
$ws.Range("K8").Value = 0.01162855693061486
$ws.Range("L8").Value = 0.01109622166105359
$ws.Range("M8").Value = 0.01087298292083813
$ws.Range("N8").Value = 0.01065861724427947
$ws.Range("O8").Value = 0.01063915752352117
$ws.Range("P8").Value = 0.01051749333642353
$ws.Range("Q8").Value = 0.01048542717561938
$ws.Range("R8").Value = 0.01047184444763432
$ws.Range("S8").Value = 0.010339418631819
$ws.Range("T8").Value = 0.01022329584523063
$ws.Range("U8").Value = 0.0102058996170457
$ws.Range("V8").Value = 0.01013531458602844
$ws.Range("W8").Value = 0.0101341529561189
$ws.Range("X8").Value = 0.01009805152874333
$ws.Range("Y8").Value = 0.01007342013260091

# Row 9
$ws.Range("C9").Value = 0.3750245571136475
$ws.Range("E9").Value = 530.4274327443891
$ws.Range("F9").Value = 0.01707904500533218
$ws.Range("G9").Value = 0.0146970752141088
$ws.Range("H9").Value = 0.01376766283675459
$ws.Range("I9").Value = 0.01313716772482882
$ws.Range("J9").Value = 0.01220680265611689
$ws.Range("K9").Value = 0.01146804324106957
$ws.Range("L9").Value = 0.01113876386502499
$ws.Range("M9").Value = 0.01113876386502499
$ws.Range("N9").Value = 0.01100592467789342
$ws.Range("O9").Value = 0.01096608183952827
$ws.Range("P9").Value = 0.01096608183952827
$ws.Range("Q9").Value = 0.0107436846064235
$ws.Range("R9").Value = 0.01051883390668531
$ws.Range("S9").Value = 0.01051081833637567
$ws.Range("T9").Value = 0.01047659330525862
$ws.Range("U9").Value = 0.01047659330525862
$ws.Range("V9").Value = 0.01041385119610058
$ws.Range("W9").Value = 0.01038207396548631
$ws.Range("X9").Value = 0.01035847482062323
$ws.Range("Y9").Value = 0.01033971603790232

# Row 10
$ws.Range("C10").Value = 0.3750009536743164
$ws.Range("E10").Value = 524.8228314972002
$ws.Range("F10").Value = 0.01715347110822386
$ws.Range("G10").Value = 0.01507617672683416
$ws.Range("H10").Value = 0.01363165984116384
$ws.Range("I10").Value = 0.01215428836823605
$ws.Range("J10").Value = 0.01156592225877968
$ws.Range("K10").Value = 0.01127611530554899
$ws.Range("L10").Value = 0.01127611530554899
$ws.Range("M10").Value = 0.01061792702703641
$ws.Range("N10").Value = 0.01061792702703641
$ws.Range("O10").Value = 0.01061792702703641
$ws.Range("P10").Value = 0.01061792702703641
$ws.Range("Q10").Value = 0.01061792702703641
$ws.Range("R10").Value = 0.01057122327037846
$ws.Range("S10").Value = 0.01046934483440373
$ws.Range("T10").Value = 0.01043869920850506
$ws.Range("U10").Value = 0.01040596253980862
$ws.Range("V10").Value = 0.01032524615299467
$ws.Range("W10").Value = 0.01030216288384276
$ws.Range("X10").Value = 0.0102606625087904
$ws.Range("Y10").Value = 0.01023046455160234

# Row 11
$ws.Range("C11").Value = 0.3593742847442627
$ws.Range("E11").Value = 509.3625587145307
$ws.Range("F11").Value = 0.01653793241387555
$ws.Range("G11").Value = 0.01451821451989983
$ws.Range("H11").Value = 0.01327107989173157
$ws.Range("I11").Value = 0.01237200521918265
$ws.Range("J11").Value = 0.01134736800426682
$ws.Range("K11").Value = 0.01134736800426682
$ws.Range("L11").Value = 0.01134736800426682
$ws.Range("M11").Value = 0.01114723044304312
$ws.Range("N11").Value = 0.01108409612104024
$ws.Range("O11").Value = 0.01042717833847545
$ws.Range("P11").Value = 0.01042717833847545
$ws.Range("Q11").Value = 0.01031665799791654
$ws.Range("R11").Value = 0.01027110679697837
$ws.Range("S11").Value = 0.01013948768262179
$ws.Range("T11").Value = 0.01013948768262179
$ws.Range("U11").Value = 0.01005659004533675
$ws.Range("V11").Value = 0.01002121577986618
$ws.Range("W11").Value = 0.009983904659773356
$ws.Range("X11").Value = 0.00996995268808535
$ws.Range("Y11").Value = 0.009929094711784221
